$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.180.98"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "2.579.54"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.23%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.580"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -8.95%  "
$ws.Range("D9").Value = "2.583.52"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.58"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.43%  "
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("D14").Value = "3.029.74"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("D15").Value = "60.191.01"
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.61%  "
$ws.Range("D18").Value = "2.594.67"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "0.0₃0841"
$ws.Range("E28").Value = "  +1.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.11"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.71%  "
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.35%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.857"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.845"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  +0.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.84"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.66%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.613"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0558"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.00%  "
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  +0.38%  "
$ws.Range("D51").Value = "1.995.17"
$ws.Range("E51").Value = "  +0.18%  "
